$wb = $excel.ActiveWorkbook

# --- Reorder sheets: climatology, grossrange, seasons (move "grossrange" before "seasons") ---
$seasons = $wb.Worksheets.Item("seasons")
$grossrange = $wb.Worksheets.Item("grossrange")
$grossrange.Move($seasons)

# re-fetch the worksheet handle by name since the sheet-position move can leave
# previously-bound variables pointing at the wrong (stale) sheet object
$grossrange = $wb.Worksheets.Item("grossrange")

# --- Update the selection on the "grossrange" sheet to the whole of column A ---
$grossrange.Columns.Item(1).Select() | Out-Null

# --- Add the new "spike" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$spike = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$spike.Name = "spike"
$spike = $wb.Worksheets.Item("spike")

# Header row (bold, matches the existing header style used elsewhere)
$spike.Range("A1").Value = "variable"
$spike.Range("B1").Value = "threshold_high"
$spike.Range("C1").Value = "threshold_low"
$spike.Range("A1:C1").Font.Bold = $true

# Data rows
$spike.Range("A2").Value = "dissolved_oxygen_uncorrected_mg_per_L"
$spike.Range("B2").Value = 5
$spike.Range("C2").Value = 2

$spike.Range("A3").Value = "dissolved_oxygen_percent_saturation"
$spike.Range("B3").Value = 10
$spike.Range("C3").Value = 5

$spike.Range("A4").Value = "salinity_psu"
$spike.Range("B4").Value = 5
$spike.Range("C4").Value = 2.5

$spike.Range("A5").Value = "sensor_depth_measured_m"
$spike.Range("B5").Value = 1
$spike.Range("C5").Value = 0.5

$spike.Range("A6").Value = "temperature_degree_C"
$spike.Range("B6").Value = 5
$spike.Range("C6").Value = 2

# Column widths (best-effort match of the authored worksheet's column widths)
$spike.Columns.Item(1).ColumnWidth = 43.584
$spike.Columns.Item(2).ColumnWidth = 13.584
$spike.Columns.Item(3).ColumnWidth = 12.417

# Select C1 and make "spike" the active sheet/tab
$spike = $wb.Worksheets.Item("spike")
$spike.Range("C1").Select() | Out-Null
$spike.Activate()

Write-Output "done"
